# Update calendar & hw
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Move the hw-0X deadline labels in column G down by one row ---
# (hw-08, hw-09, hw-10, hw-11 each shift from the "lecture" row to the
#  following row)
$moves = @(
    @{ From = "G26"; To = "G27" },
    @{ From = "G29"; To = "G30" },
    @{ From = "G32"; To = "G33" },
    @{ From = "G35"; To = "G36" }
)

foreach ($m in $moves) {
    $val = $ws.Range($m.From).Value2
    $ws.Range($m.From).Value = $null
    $ws.Range($m.To).Value = $val
}

# --- New hw-12 deadline label ---
$ws.Range("G39").Value = "hw-12"

# --- Insert a new row above row 40 for the "Fri, Nov 22" entry ---
$ws.Rows.Item(40).Insert()
$ws.Rows.Item(40).RowHeight = 18

$ws.Range("B40").Value = "Fri, Nov 22"

# --- Update the active selection shown in the workbook ---
$ws.Range("B40").Select() | Out-Null
